# ------------------------------------------------------------------
# Applies the commit:
#  - Rename sheet "Phạt" -> "Đơn sale phụ"
#  - Populate "Đơn sale chính" (sheet1) with a sale order row + total row
#  - Replace contents of "Đơn sale phụ" (sheet2, was "Phạt") with a
#    sale-order-style header/row/total (same columns as sheet1 but for
#    sale phụ)
#  - Update the "Lương" (sheet3) summary sheet: delete the old
#    "Phạt tại LONG XUYÊN" row and refresh several computed totals
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # "Đơn sale chính"
$ws2 = $wb.Worksheets.Item(2)   # "Phạt" -> "Đơn sale phụ"
$ws3 = $wb.Worksheets.Item(3)   # "Lương"

# ------------------------------------------------------------------
# 1. Rename sheet 2
# ------------------------------------------------------------------
$ws2.Name = "Đơn sale phụ"

# ------------------------------------------------------------------
# 2. "Đơn sale chính" sheet: write headers + data row + totals row
# ------------------------------------------------------------------
$headers1 = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale chính","Chiết khấu sale chính")
for ($i = 0; $i -lt $headers1.Length; $i++) {
    $ws1.Cells.Item(1, $i + 1).Value = $headers1[$i]
}

$ws1.Cells.Item(2, 1).Value = "HD-LUXURY"
$ws1.Cells.Item(2, 2).Value = 633
$ws1.Cells.Item(2, 3).NumberFormat = "@"
$ws1.Cells.Item(2, 3).Value = "08-07-2024"
$ws1.Cells.Item(2, 4).Value = "CẦN THƠ"
$ws1.Cells.Item(2, 5).Value = "Bạch Nhi"
$ws1.Cells.Item(2, 6).Value = "Cá nhân"
$ws1.Cells.Item(2, 7).Value = "Nâng mũi"
$ws1.Cells.Item(2, 8).Value = 15000000
# I2, J2 left blank (Sale phụ / Upsale not applicable on this row)
$ws1.Cells.Item(2, 11).Value = 15000000
$ws1.Cells.Item(2, 12).Value = 7000000
$ws1.Cells.Item(2, 13).Value = 0.1
$ws1.Cells.Item(2, 14).Value = 700000

$ws1.Cells.Item(3, 1).Value = "Tổng"
$ws1.Cells.Item(3, 2).Value = 1
$ws1.Cells.Item(3, 8).Value = 15000000
$ws1.Cells.Item(3, 10).Value = 0
$ws1.Cells.Item(3, 11).Value = 15000000
$ws1.Cells.Item(3, 12).Value = 7000000
$ws1.Cells.Item(3, 13).Value = 0
$ws1.Cells.Item(3, 14).Value = 700000

# ------------------------------------------------------------------
# 3. "Đơn sale phụ" sheet (previously "Phạt"): clear old contents and
#    write the new headers + data row + totals row
# ------------------------------------------------------------------
$ws2.Cells.Clear()

$headers2 = @("Tiền tố","Mã dịch vụ","Ngày thực hiện","Cơ sở","Khách hàng","Nguồn khách","Tên dịch vụ","Đơn giá gốc","Sale phụ","Upsale","Đơn giá","Đã thanh toán","Tỉ lệ chiết khấu sale phụ","Chiết khấu sale phụ")
for ($i = 0; $i -lt $headers2.Length; $i++) {
    $ws2.Cells.Item(1, $i + 1).Value = $headers2[$i]
}

$ws2.Cells.Item(2, 1).Value = "HD-LUXURY"
$ws2.Cells.Item(2, 2).Value = 625
$ws2.Cells.Item(2, 3).NumberFormat = "@"
$ws2.Cells.Item(2, 3).Value = "08-04-2024"
$ws2.Cells.Item(2, 4).Value = "SÓC TRĂNG"
$ws2.Cells.Item(2, 5).Value = "nguyễn thị mỹ chăm"
$ws2.Cells.Item(2, 6).Value = "Cá nhân"
$ws2.Cells.Item(2, 7).Value = "Cắt mí"
# H2 (Đơn giá gốc) left blank
$ws2.Cells.Item(2, 9).Value = "Lê Hoàng Thanh"
$ws2.Cells.Item(2, 10).Value = 6000000
$ws2.Cells.Item(2, 11).Value = 6000000
$ws2.Cells.Item(2, 12).Value = 6000000
$ws2.Cells.Item(2, 13).Value = 0.04
$ws2.Cells.Item(2, 14).Value = 0

$ws2.Cells.Item(3, 1).Value = "Tổng"
$ws2.Cells.Item(3, 2).Value = 1
$ws2.Cells.Item(3, 8).Value = 0
$ws2.Cells.Item(3, 10).Value = 6000000
$ws2.Cells.Item(3, 11).Value = 6000000
$ws2.Cells.Item(3, 12).Value = 6000000
$ws2.Cells.Item(3, 13).Value = 0
$ws2.Cells.Item(3, 14).Value = 0

# ------------------------------------------------------------------
# 4. "Lương" sheet: remove the obsolete "Phạt tại LONG XUYÊN" row
#    (row 21), which shifts every following "tại SÓC TRĂNG" / totals
#    row up by one, then refresh the recomputed totals.
# ------------------------------------------------------------------
$ws3.Rows.Item(21).Delete()

$ws3.Cells.Item(5, 2).Value  = 700000              # Chiết khấu sale chính tại CẦN THƠ
$ws3.Cells.Item(12, 2).Value = 8.5                 # Tổng công tại LONG XUYÊN
$ws3.Cells.Item(13, 2).Value = 1517857.142857143   # Lương cơ bản tại LONG XUYÊN
$ws3.Cells.Item(23, 2).ClearContents()             # Lương cơ bản tại SÓC TRĂNG (blank)
$ws3.Cells.Item(31, 2).Value = 700000              # Tổng lương tại CẦN THƠ
$ws3.Cells.Item(32, 2).Value = 1517857.142857143   # Tổng lương tại LONG XUYÊN
$ws3.Cells.Item(34, 2).Value = 2217857.142857143   # Tổng lương tại HỆ THỐNG
